# Feature_Tracker.xlsx - "Added a new feature to the feature tracker."
#
# 1) Correct the Requestor for the existing "Make average more accurate"
#    row (row 26): it should read "Stefan Titus" instead of "Weston Fiala".
# 2) Append a brand-new feature request row (row 38):
#       Feature     : Allow for different dice images
#       Description : I want to be able to change the images that my dice use.
#       Requestors  : Rae - urbanchika@gmail.com
#    (Completed Version column D left blank, same as other un-shipped rows.)
# 3) Move the active-cell selection down to B39, matching where the user's
#    cursor ends up after typing the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the Requestor on the "Make average more accurate" row ---
$ws.Range("D26").Value = "Stefan Titus"

# --- 2. Add the new feature request as row 38 ---
$ws.Range("A38").Value = "Allow for different dice images"
$ws.Range("B38").Value = "I want to be able to change the images that my dice use."
$ws.Range("D38").Value = "Rae - urbanchika@gmail.com"

# --- 3. Update the selection to reflect the next empty row ---
$ws.Range("B39").Select() | Out-Null
